$d = $word.ActiveDocument

# Locate the paragraph that ends with " vs LW 1,11 %" (the "Lack of Couriers" line)
# so we can insert the new "level 3" paragraph immediately after it.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*vs LW 1,11 %*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# Split off a new paragraph right after the target paragraph; this keeps the
# target paragraph's own formatting/spacing untouched and gives the new
# paragraph the same "w:spacing w:after=0" paragraph formatting (inherited
# from the paragraph it was split from).
$target.Range.InsertParagraphAfter()

# Re-fetch the freshly created (now existing) paragraph by its index and set
# its text.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Dentro de level 3 se identifico que la razón wo es automation_lack_of_rts_other con 0,73 % vs LW 0,77 %"
